# Update leve-profit tracking figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to reflect refreshed market-board pricing pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1058.9454
$ws.Range("I15").Value = 1058.9454
$ws.Range("K15").Value = 3176.8362
$ws.Range("M15").Value = -3007.8362
$ws.Range("H47").Value = 341766.5
$ws.Range("I47").Value = 444355.34
$ws.Range("J47").Value = 34000
$ws.Range("K47").Value = 444355.34
$ws.Range("L47").Value = 34000
$ws.Range("M47").Value = -443383.34
$ws.Range("N47").Value = -35944
$ws.Range("H97").Value = 2809.625
$ws.Range("J97").Value = 2809.625
$ws.Range("L97").Value = 8428.875
$ws.Range("N97").Value = -9420.875
$ws.Range("H100").Value = 5367
$ws.Range("I100").Value = 3873.818
$ws.Range("K100").Value = 3873.818
$ws.Range("M100").Value = -3332.818
$ws.Range("H106").Value = 8732.214
$ws.Range("I106").Value = 9772
$ws.Range("J106").Value = 4919.6665
$ws.Range("K106").Value = 9772
$ws.Range("L106").Value = 4919.6665
$ws.Range("M106").Value = -9141
$ws.Range("N106").Value = -6181.6665
$ws.Range("H132").Value = 2411.348
$ws.Range("I132").Value = 1987.579
$ws.Range("J132").Value = 4424.25
$ws.Range("K132").Value = 5962.737
$ws.Range("L132").Value = 13272.75
$ws.Range("M132").Value = -3432.737
$ws.Range("N132").Value = -18332.75
$ws.Range("H141").Value = 5702.6875
$ws.Range("I141").Value = 6187.3335
$ws.Range("K141").Value = 18562.0005
$ws.Range("M141").Value = -13382.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 500050000
$ws.Range("J17").Value = 500050000
$ws.Range("L17").Value = 500050000
$ws.Range("N17").Value = -500050346
$ws.Range("H32").Value = 10528.85
$ws.Range("I32").Value = 9744.200000000001
$ws.Range("K32").Value = 9744.200000000001
$ws.Range("M32").Value = -9457.200000000001
$ws.Range("H61").Value = 14204521
$ws.Range("I61").Value = 16670896
$ws.Range("J61").Value = 1256051.5
$ws.Range("K61").Value = 16670896
$ws.Range("L61").Value = 1256051.5
$ws.Range("M61").Value = -16670684
$ws.Range("N61").Value = -1256475.5
$ws.Range("H74").Value = 1812.1333
$ws.Range("I74").Value = 1847.7142
$ws.Range("J74").Value = 1314
$ws.Range("K74").Value = 1847.7142
$ws.Range("L74").Value = 1314
$ws.Range("M74").Value = -973.7141999999999
$ws.Range("N74").Value = -3062
$ws.Range("H77").Value = 1812.1333
$ws.Range("I77").Value = 1847.7142
$ws.Range("J77").Value = 1314
$ws.Range("K77").Value = 9238.571
$ws.Range("L77").Value = 6570
$ws.Range("M77").Value = -4870.571
$ws.Range("N77").Value = -15306
$ws.Range("H122").Value = 3624.625
$ws.Range("I122").Value = 4010
$ws.Range("J122").Value = 2468.5
$ws.Range("K122").Value = 12030
$ws.Range("L122").Value = 7405.5
$ws.Range("M122").Value = -9580
$ws.Range("N122").Value = -12305.5
$ws.Range("H136").Value = 14204521
$ws.Range("I136").Value = 16670896
$ws.Range("J136").Value = 1256051.5
$ws.Range("K136").Value = 50012688
$ws.Range("L136").Value = 3768154.5
$ws.Range("M136").Value = -50010138
$ws.Range("N136").Value = -3773254.5
$ws.Range("H139").Value = 149285.58
$ws.Range("J139").Value = 149285.58
$ws.Range("L139").Value = 149285.58
$ws.Range("N139").Value = -159565.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2219
$ws.Range("I94").Value = 2478.2173
$ws.Range("K94").Value = 2478.2173
$ws.Range("M94").Value = -2027.2173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62.52174
$ws.Range("J7").Value = 56.9
$ws.Range("L7").Value = 56.9
$ws.Range("N7").Value = -282.9
$ws.Range("H22").Value = 312.7143
$ws.Range("I22").Value = 314.83334
$ws.Range("K22").Value = 314.83334
$ws.Range("M22").Value = 35.16665999999998
$ws.Range("H31").Value = 20836162
$ws.Range("I31").Value = 30305378
$ws.Range("K31").Value = 30305378
$ws.Range("M31").Value = -30305083
$ws.Range("H32").Value = 500005000
$ws.Range("I32").Value = 500005000
$ws.Range("K32").Value = 500005000
$ws.Range("M32").Value = -500004684
$ws.Range("H34").Value = 20836162
$ws.Range("I34").Value = 30305378
$ws.Range("K34").Value = 30305378
$ws.Range("M34").Value = -30305176
$ws.Range("H35").Value = 18000
$ws.Range("I35").Value = 18000
$ws.Range("K35").Value = 18000
$ws.Range("M35").Value = -17706
$ws.Range("H131").Value = 99998
$ws.Range("J131").Value = 99998
$ws.Range("L131").Value = 99998
$ws.Range("N131").Value = -110078
$ws.Range("H132").Value = 2080.8572
$ws.Range("I132").Value = 2347.3125
$ws.Range("K132").Value = 7041.9375
$ws.Range("M132").Value = -4511.9375
$ws.Range("H134").Value = 2036.7941
$ws.Range("I134").Value = 1813.1154
$ws.Range("K134").Value = 5439.3462
$ws.Range("M134").Value = -2904.3462
$ws.Range("H141").Value = 655865.2
$ws.Range("I141").Value = 200000
$ws.Range("K141").Value = 200000
$ws.Range("M141").Value = -194820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 6352.5
$ws.Range("J60").Value = 18408.25
$ws.Range("L60").Value = 55224.75
$ws.Range("N60").Value = -55726.75
$ws.Range("H124").Value = 8070.1113
$ws.Range("I124").Value = 1766.6666
$ws.Range("J124").Value = 11221.833
$ws.Range("K124").Value = 5299.9998
$ws.Range("L124").Value = 33665.499
$ws.Range("M124").Value = -389.9997999999996
$ws.Range("N124").Value = -43485.499
$ws.Range("H131").Value = 4336.7393
$ws.Range("I131").Value = 2080.4546
$ws.Range("K131").Value = 6241.3638
$ws.Range("M131").Value = -1201.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 97326
$ws.Range("J45").Value = 97326
$ws.Range("L45").Value = 97326
$ws.Range("N45").Value = -98444
$ws.Range("H122").Value = 7964.5713
$ws.Range("I122").Value = 6047.3125
$ws.Range("K122").Value = 18141.9375
$ws.Range("M122").Value = -15691.9375
$ws.Range("H132").Value = 2582730.5
$ws.Range("I132").Value = 2930.2903
$ws.Range("J132").Value = 15911698
$ws.Range("K132").Value = 8790.8709
$ws.Range("L132").Value = 47735094
$ws.Range("M132").Value = -6260.8709
$ws.Range("N132").Value = -47740154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7968.3
$ws.Range("I22").Value = 17340.715
$ws.Range("J22").Value = 2921.6155
$ws.Range("K22").Value = 17340.715
$ws.Range("L22").Value = 2921.6155
$ws.Range("M22").Value = -17045.715
$ws.Range("N22").Value = -3511.6155
$ws.Range("H27").Value = 7968.3
$ws.Range("I27").Value = 17340.715
$ws.Range("J27").Value = 2921.6155
$ws.Range("K27").Value = 17340.715
$ws.Range("L27").Value = 2921.6155
$ws.Range("M27").Value = -17233.715
$ws.Range("N27").Value = -3135.6155
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H46").Value = 997.05884
$ws.Range("J46").Value = 1912
$ws.Range("L46").Value = 1912
$ws.Range("N46").Value = -2288
$ws.Range("H55").Value = 1094.7407
$ws.Range("I55").Value = 791
$ws.Range("J55").Value = 1421.8462
$ws.Range("K55").Value = 791
$ws.Range("L55").Value = 1421.8462
$ws.Range("M55").Value = -618
$ws.Range("N55").Value = -1767.8462
$ws.Range("H82").Value = 2799.2
$ws.Range("I82").Value = 1017.6
$ws.Range("J82").Value = 4580.8
$ws.Range("K82").Value = 1017.6
$ws.Range("L82").Value = 4580.8
$ws.Range("M82").Value = -656.6
$ws.Range("N82").Value = -5302.8
$ws.Range("H85").Value = 2799.2
$ws.Range("I85").Value = 1017.6
$ws.Range("J85").Value = 4580.8
$ws.Range("K85").Value = 1017.6
$ws.Range("L85").Value = 4580.8
$ws.Range("M85").Value = 230.4
$ws.Range("N85").Value = -7076.8
$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079
$ws.Range("H136").Value = 5382.1333
$ws.Range("I136").Value = 2456
$ws.Range("J136").Value = 7332.8887
$ws.Range("K136").Value = 7368
$ws.Range("L136").Value = 21998.6661
$ws.Range("M136").Value = -4818
$ws.Range("N136").Value = -27098.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H47").Value = 39429.75
$ws.Range("J47").Value = 39429.75
$ws.Range("L47").Value = 39429.75
$ws.Range("N47").Value = -40573.75
$ws.Range("H48").Value = 40062.668
$ws.Range("I48").Value = 40059
$ws.Range("J48").Value = 40064.5
$ws.Range("K48").Value = 40059
$ws.Range("L48").Value = 40064.5
$ws.Range("M48").Value = -39490
$ws.Range("N48").Value = -41202.5
$ws.Range("H100").Value = 1019.15625
$ws.Range("I100").Value = 987.087
$ws.Range("K100").Value = 1974.174
$ws.Range("M100").Value = -1433.174
$ws.Range("H132").Value = 419986.34
$ws.Range("I132").Value = 3496.1082
$ws.Range("J132").Value = 1820908
$ws.Range("K132").Value = 10488.3246
$ws.Range("L132").Value = 5462724
$ws.Range("M132").Value = -7958.3246
$ws.Range("N132").Value = -5467784
$ws.Range("H136").Value = 591429.4399999999
$ws.Range("I136").Value = 3905.5833
$ws.Range("J136").Value = 2001486.6
$ws.Range("K136").Value = 11716.7499
$ws.Range("L136").Value = 6004459.800000001
$ws.Range("M136").Value = -9166.749899999999
$ws.Range("N136").Value = -6009559.800000001

